$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade record (row 4) with the same layout as the existing
# rows: Date, Profitable, Principle, Start Principle, BuyPrice, SellPrice,
# IsShortSell, Price Change %, Strong trade.
$ws.Cells.Item(4, 1).Value = 42635.800439814811
$ws.Cells.Item(4, 2).Value = $true
$ws.Cells.Item(4, 3).Value = 10085.129999999999
$ws.Cells.Item(4, 4).Value = 10020
$ws.Cells.Item(4, 5).Value = 82.03
$ws.Cells.Item(4, 6).Value = 80.9599
$ws.Cells.Item(4, 7).Value = $true
$ws.Cells.Item(4, 8).Value = -1.3
$ws.Cells.Item(4, 9).Value = $false

# Match the date-style formatting used for the Date and IsShortSell columns
# in the preceding row (row 3) by copying their formats down to row 4.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
